$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$shp = $nm.Shapes.Item(2)
Write-Output ("before=" + $shp.TextFrame.TextRange.Text)
try {
    $shp.TextFrame.TextRange.Text = "HELLO"
    Write-Output ("immediately after=" + $shp.TextFrame.TextRange.Text)
} catch {
    Write-Output ("EXC: " + $_)
}
$nm2 = $p.NotesMaster
Write-Output ("reload nm shape=" + $nm2.Shapes.Item(2).TextFrame.TextRange.Text)
